# Scheduled-runner price refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H, I, J, K, L, M, N) across the leve-profit sheets to
# reflect newly pulled market data. Mirrors the upstream diff cell-by-cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 37557
$ws.Range("J123").Value = 37557
$ws.Range("L123").Value = 37557
$ws.Range("N123").Value = -47357
$ws.Range("H124").Value = 43706
$ws.Range("J124").Value = 43706
$ws.Range("L124").Value = 43706
$ws.Range("N124").Value = -53526
$ws.Range("H126").Value = 42305.6
$ws.Range("J126").Value = 42305.6
$ws.Range("L126").Value = 42305.6
$ws.Range("N126").Value = -52185.6
$ws.Range("H128").Value = 42516.25
$ws.Range("J128").Value = 42516.25
$ws.Range("L128").Value = 42516.25
$ws.Range("N128").Value = -52476.25
$ws.Range("H130").Value = 43114.668
$ws.Range("J130").Value = 43114.668
$ws.Range("L130").Value = 43114.668
$ws.Range("N130").Value = -53154.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13193.38
$ws.Range("I32").Value = 13057.293
$ws.Range("K32").Value = 13057.293
$ws.Range("M32").Value = -12770.293
$ws.Range("H109").Value = 36856.75
$ws.Range("J109").Value = 36856.75
$ws.Range("L109").Value = 36856.75
$ws.Range("N109").Value = -39630.75
$ws.Range("H111").Value = 49451
$ws.Range("J111").Value = 49451
$ws.Range("L111").Value = 49451
$ws.Range("N111").Value = -57631
$ws.Range("H117").Value = 49561.75
$ws.Range("J117").Value = 49561.75
$ws.Range("L117").Value = 49561.75
$ws.Range("N117").Value = -58739.75
$ws.Range("H118").Value = 49803
$ws.Range("J118").Value = 49803
$ws.Range("L118").Value = 49803
$ws.Range("N118").Value = -53117
$ws.Range("H119").Value = 27658.4
$ws.Range("J119").Value = 27658.4
$ws.Range("L119").Value = 27658.4
$ws.Range("N119").Value = -37334.4
$ws.Range("H121").Value = 38248.2
$ws.Range("J121").Value = 38248.2
$ws.Range("L121").Value = 38248.2
$ws.Range("N121").Value = -41742.2
$ws.Range("H122").Value = 1934.1428
$ws.Range("I122").Value = 1927.2106
$ws.Range("K122").Value = 5781.6318
$ws.Range("M122").Value = -3331.6318
$ws.Range("H123").Value = 47936.332
$ws.Range("J123").Value = 47936.332
$ws.Range("L123").Value = 47936.332
$ws.Range("N123").Value = -57736.332
$ws.Range("H125").Value = 32698.143
$ws.Range("J125").Value = 32698.143
$ws.Range("L125").Value = 32698.143
$ws.Range("N125").Value = -42538.143
$ws.Range("H130").Value = 38462.332
$ws.Range("J130").Value = 38462.332
$ws.Range("L130").Value = 38462.332
$ws.Range("N130").Value = -48502.332
$ws.Range("H132").Value = 9435579
$ws.Range("I132").Value = 15626145
$ws.Range("J132").Value = 2336.8572
$ws.Range("K132").Value = 46878435
$ws.Range("L132").Value = 7010.571599999999
$ws.Range("M132").Value = -46875905
$ws.Range("N132").Value = -12070.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 46838
$ws.Range("J108").Value = 46838
$ws.Range("L108").Value = 46838
$ws.Range("N108").Value = -54518
$ws.Range("H117").Value = 47498
$ws.Range("J117").Value = 47498
$ws.Range("L117").Value = 47498
$ws.Range("N117").Value = -56676
$ws.Range("H122").Value = 40722.6
$ws.Range("J122").Value = 40722.6
$ws.Range("L122").Value = 40722.6
$ws.Range("N122").Value = -50522.6
$ws.Range("H130").Value = 49181
$ws.Range("J130").Value = 49181
$ws.Range("L130").Value = 49181
$ws.Range("N130").Value = -59221

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49885.5
$ws.Range("J20").Value = 49885.5
$ws.Range("L20").Value = 49885.5
$ws.Range("N20").Value = -50357.5
$ws.Range("H30").Value = 49885.5
$ws.Range("J30").Value = 49885.5
$ws.Range("L30").Value = 49885.5
$ws.Range("N30").Value = -50067.5
$ws.Range("H100").Value = 44617.332
$ws.Range("J100").Value = 44617.332
$ws.Range("L100").Value = 44617.332
$ws.Range("N100").Value = -46781.332
$ws.Range("H128").Value = 49885.5
$ws.Range("J128").Value = 49885.5
$ws.Range("L128").Value = 49885.5
$ws.Range("N128").Value = -59845.5
$ws.Range("H133").Value = 27443.2
$ws.Range("J133").Value = 27443.2
$ws.Range("L133").Value = 27443.2
$ws.Range("N133").Value = -32503.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4441.5
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 5004.5713
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 15013.7139
$ws.Range("M118").Value = -257
$ws.Range("N118").Value = -17499.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 46701
$ws.Range("J110").Value = 46701
$ws.Range("L110").Value = 46701
$ws.Range("N110").Value = -54881
$ws.Range("H114").Value = 45628.25
$ws.Range("J114").Value = 45628.25
$ws.Range("L114").Value = 45628.25
$ws.Range("N114").Value = -54306.25
$ws.Range("H116").Value = 39000
$ws.Range("J116").Value = 39000
$ws.Range("L116").Value = 39000
$ws.Range("N116").Value = -48178
$ws.Range("H122").Value = 1243.8572
$ws.Range("I122").Value = 1002.3333
$ws.Range("J122").Value = 1425
$ws.Range("K122").Value = 3006.9999
$ws.Range("L122").Value = 4275
$ws.Range("M122").Value = -556.9998999999998
$ws.Range("N122").Value = -9175
$ws.Range("H124").Value = 38420.332
$ws.Range("J124").Value = 41776
$ws.Range("L124").Value = 41776
$ws.Range("N124").Value = -51596
$ws.Range("H128").Value = 38996
$ws.Range("J128").Value = 38996
$ws.Range("L128").Value = 38996
$ws.Range("N128").Value = -48956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2930.6667
$ws.Range("I7").Value = 2158.6875
$ws.Range("J7").Value = 5401
$ws.Range("K7").Value = 2158.6875
$ws.Range("L7").Value = 5401
$ws.Range("M7").Value = -2046.6875
$ws.Range("N7").Value = -5625
$ws.Range("H36").Value = 48696
$ws.Range("J36").Value = 48696
$ws.Range("L36").Value = 48696
$ws.Range("N36").Value = -49820
$ws.Range("H40").Value = 5648.3
$ws.Range("I40").Value = 3116.3333
$ws.Range("J40").Value = 9446.25
$ws.Range("K40").Value = 3116.3333
$ws.Range("L40").Value = 9446.25
$ws.Range("M40").Value = -2980.3333
$ws.Range("N40").Value = -9718.25
$ws.Range("H122").Value = 93010.17999999999
$ws.Range("I122").Value = 102111.2
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 306333.6
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -303883.6
$ws.Range("N122").Value = -10900
$ws.Range("H123").Value = 39425
$ws.Range("J123").Value = 39425
$ws.Range("L123").Value = 39425
$ws.Range("N123").Value = -49225
$ws.Range("H125").Value = 48211
$ws.Range("J125").Value = 48211
$ws.Range("L125").Value = 48211
$ws.Range("N125").Value = -58051
$ws.Range("H126").Value = 2930.6667
$ws.Range("I126").Value = 2158.6875
$ws.Range("J126").Value = 5401
$ws.Range("K126").Value = 6476.0625
$ws.Range("L126").Value = 16203
$ws.Range("M126").Value = -4006.0625
$ws.Range("N126").Value = -21143
$ws.Range("H127").Value = 50602.75
$ws.Range("J127").Value = 50602.75
$ws.Range("L127").Value = 50602.75
$ws.Range("N127").Value = -60522.75
$ws.Range("H128").Value = 34425
$ws.Range("J128").Value = 34425
$ws.Range("L128").Value = 34425
$ws.Range("N128").Value = -44385

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 42284
$ws.Range("J121").Value = 42284
$ws.Range("L121").Value = 42284
$ws.Range("N121").Value = -45778
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
